$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (143-145), columns A (1) through AM (39)
$data = @(
    @(1481,26439,8585,715,1228,1181,201801,2018,1,1813,1895,2077,1761,1641,24637,26265,23009,24107,21732,8181,10300,8115,6582,5050,709,1077,840,1356,1121,1441,1548,1386,2138,1155,1488,1061,886,815,991),
    @(1307,23976,6749,857,1137,1150,201802,2018,2,1481,1813,1895,2152,1216,26439,24637,26265,25699,16549,8585,8181,10300,6298,4379,715,709,1077,1062,570,1228,1441,1548,1685,1020,1181,1488,1061,2015,872),
    @(1403,29022,8796,1322,1834,2528,201803,2018,3,1307,1481,1813,2187,1808,23976,26439,24637,26317,27056,6749,8585,8181,8613,4938,857,715,709,959,1117,1137,1228,1441,1654,1855,1150,1181,1488,1132,1405)
)

$startRow = 143
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $data[$i]
    for ($col = 1; $col -le $rowValues.Length; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $rowValues[$col - 1]
    }
}
